$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.253.63'
$ws.Range("E2").Value = '  -3.11%  '
$ws.Range("D3").Value = '2.462.16'
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.00'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.02'
$ws.Range("E6").Value = '  -6.49%  '
$ws.Range("E7").Value = '  -3.30%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -4.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.45'
$ws.Range("E10").Value = '  -6.95%  '
$ws.Range("E11").Value = '  -2.60%  '
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.92'
$ws.Range("E13").Value = '  -4.89%  '
$ws.Range("D14").Value = '2.837.75'
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '2.443.08'
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.61'
$ws.Range("E16").Value = '  -6.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("D18").Value = '41.203.89'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("E19").Value = '  -5.61%  '
$ws.Range("D20").Value = '0.0₃0918'
$ws.Range("E20").Value = '  -2.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.44'
$ws.Range("E21").Value = '  -5.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.78'
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.99'
$ws.Range("E23").Value = '  -3.03%  '
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("E25").Value = '  -5.46%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.32'
$ws.Range("E27").Value = '  -5.88%  '
$ws.Range("E28").Value = '  -4.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("E30").Value = '  -7.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.67'
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.54'
$ws.Range("E32").Value = '  -3.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.59'
$ws.Range("E33").Value = '  -6.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("E35").Value = '  -4.28%  '
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("E37").Value = '  -6.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.01'
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("E40").Value = '  -8.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.23'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.90'
$ws.Range("E42").Value = '  -5.60%  '
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").Value = '1.960.74'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("E45").Value = '  -4.61%  '
$ws.Range("E46").Value = '  -6.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.64'
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '69.79'
$ws.Range("E48").Value = '  -3.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.72'
$ws.Range("E49").Value = '  -5.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.14'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  -5.91%  '
